# Backup of the correct code:
#  - B2 on the first worksheet should hold "Sauce Labs Backpack" instead of
#    "Sauce Labs Onesie" (the now-unused "Sauce Labs Onesie" shared string
#    is dropped automatically once nothing references it anymore).
#  - The active selection on that worksheet moves to B6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "Sauce Labs Backpack"

$ws.Range("B6").Select() | Out-Null
